$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force column D to Text format so numeric-looking strings
# (e.g. "0.08452", "0.00001112") are stored as literal text instead of
# being parsed/rounded as floating point numbers by Excel.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.226.40"
$ws.Range("D3").Value = "1.897.99"
$ws.Range("D5").Value = "325.36"
$ws.Range("D7").Value = "0.5160"
$ws.Range("D8").Value = "0.4009"
$ws.Range("D9").Value = "0.08452"
$ws.Range("D10").Value = "42.67"
$ws.Range("D11").Value = "1.118"
$ws.Range("D12").Value = "23.36"
$ws.Range("D14").Value = "1.892.84"
$ws.Range("D15").Value = "7.337"
$ws.Range("D17").Value = "94.79"
$ws.Range("D18").Value = "0.00001112"
$ws.Range("D19").Value = "0.06648"
$ws.Range("D20").Value = "18.24"
$ws.Range("D23").Value = "30.227.29"
$ws.Range("D24").Value = "11.28"
$ws.Range("D26").Value = "2.110.23"
$ws.Range("D27").Value = "21.63"
$ws.Range("D28").Value = "161.17"
$ws.Range("D29").Value = "2.361"
$ws.Range("D30").Value = "129.07"
$ws.Range("D31").Value = "1.097"
$ws.Range("D32").Value = "0.1058"
$ws.Range("D33").Value = "6.074"
$ws.Range("D35").Value = "0.02498"
$ws.Range("D36").Value = "0.06560"
$ws.Range("D37").Value = "5.257"
$ws.Range("D38").Value = "0.2201"
$ws.Range("D40").Value = "11.80"
$ws.Range("D42").Value = "8.727"
$ws.Range("D43").Value = "1.232"
$ws.Range("D44").Value = "0.6103"
$ws.Range("D45").Value = "13.27"
$ws.Range("D47").Value = "2.057"
$ws.Range("D49").Value = "124.44"
$ws.Range("D50").Value = "1.166"
$ws.Range("D51").Value = "78.99"

# Restore original (unstyled) formatting on column D now that the
# literal text values have been committed.
$dRange.ClearFormats()

# Column E values are percentage strings with surrounding whitespace
# and are never auto-converted to numbers, so no special handling is
# required for them.
$ws.Range("E2").Value = "  +3.14%  "
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("E5").Value = "  +3.45%  "
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +2.03%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +12.94%  "
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +3.62%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  -2.86%  "
$ws.Range("E30").Value = "  +1.21%  "
$ws.Range("E31").Value = "  +3.41%  "
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("E35").Value = "  +0.55%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("E40").Value = "  +4.77%  "
$ws.Range("E41").Value = "  -0.30%  "
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  +1.83%  "
